# Auto-generated edit script: apply numeric updates per diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Cells.Item(19, 8).Value = 7655.2144
$ws.Cells.Item(19, 9).Value = 1696.8125
$ws.Cells.Item(19, 10).Value = 15599.75
$ws.Cells.Item(19, 11).Value = 1696.8125
$ws.Cells.Item(19, 12).Value = 15599.75
$ws.Cells.Item(19, 13).Value = -1521.8125
$ws.Cells.Item(19, 14).Value = -15949.75
$ws.Cells.Item(88, 8).Value = 17550782
$ws.Cells.Item(88, 9).Value = 37038520
$ws.Cells.Item(88, 10).Value = 11818.7
$ws.Cells.Item(88, 11).Value = 37038520
$ws.Cells.Item(88, 12).Value = 11818.7
$ws.Cells.Item(88, 13).Value = -37038114
$ws.Cells.Item(88, 14).Value = -12630.7
$ws.Cells.Item(91, 8).Value = 17550782
$ws.Cells.Item(91, 9).Value = 37038520
$ws.Cells.Item(91, 10).Value = 11818.7
$ws.Cells.Item(91, 11).Value = 37038520
$ws.Cells.Item(91, 12).Value = 11818.7
$ws.Cells.Item(91, 13).Value = -37037116
$ws.Cells.Item(91, 14).Value = -14626.7
$ws.Cells.Item(132, 8).Value = 4290.705
$ws.Cells.Item(132, 9).Value = 2102.64
$ws.Cells.Item(132, 10).Value = 14236.454
$ws.Cells.Item(132, 11).Value = 6307.92
$ws.Cells.Item(132, 12).Value = 42709.362
$ws.Cells.Item(132, 13).Value = -3777.92
$ws.Cells.Item(132, 14).Value = -47769.362
$ws.Cells.Item(137, 8).Value = 4879.647
$ws.Cells.Item(137, 9).Value = 5158.878
$ws.Cells.Item(137, 10).Value = 3734.8
$ws.Cells.Item(137, 11).Value = 15476.634
$ws.Cells.Item(137, 12).Value = 11204.4
$ws.Cells.Item(137, 13).Value = -12926.634
$ws.Cells.Item(137, 14).Value = -16304.4

$ws = $wb.Worksheets("ARM")
$ws.Cells.Item(74, 8).Value = 1104.5
$ws.Cells.Item(74, 9).Value = 1104.5
$ws.Cells.Item(74, 11).Value = 1104.5
$ws.Cells.Item(74, 13).Value = -230.5
$ws.Cells.Item(77, 8).Value = 1104.5
$ws.Cells.Item(77, 9).Value = 1104.5
$ws.Cells.Item(77, 11).Value = 5522.5
$ws.Cells.Item(77, 13).Value = -1154.5

$ws = $wb.Worksheets("BSM")
$ws.Cells.Item(68, 8).Value = 49999
$ws.Cells.Item(68, 10).Value = 49999
$ws.Cells.Item(68, 12).Value = 49999
$ws.Cells.Item(68, 14).Value = -51621
$ws.Cells.Item(71, 8).Value = 49999
$ws.Cells.Item(71, 10).Value = 49999
$ws.Cells.Item(71, 12).Value = 149997
$ws.Cells.Item(71, 14).Value = -158109
$ws.Cells.Item(75, 8).Value = 11495
$ws.Cells.Item(75, 9).Value = 2990
$ws.Cells.Item(75, 11).Value = 2990
$ws.Cells.Item(75, 13).Value = -2054
$ws.Cells.Item(78, 8).Value = 11495
$ws.Cells.Item(78, 9).Value = 2990
$ws.Cells.Item(78, 11).Value = 8970
$ws.Cells.Item(78, 13).Value = -4290
$ws.Cells.Item(82, 8).Value = 17900
$ws.Cells.Item(82, 9).Value = 4737.75
$ws.Cells.Item(82, 11).Value = 4737.75
$ws.Cells.Item(82, 13).Value = -4354.75
$ws.Cells.Item(85, 8).Value = 17900
$ws.Cells.Item(85, 9).Value = 4737.75
$ws.Cells.Item(85, 11).Value = 4737.75
$ws.Cells.Item(85, 13).Value = -3411.75
$ws.Cells.Item(86, 8).Value = 1983.875
$ws.Cells.Item(86, 9).Value = 1221.5
$ws.Cells.Item(86, 10).Value = 2746.25
$ws.Cells.Item(86, 11).Value = 1221.5
$ws.Cells.Item(86, 12).Value = 2746.25
$ws.Cells.Item(86, 13).Value = -98.5
$ws.Cells.Item(86, 14).Value = -4992.25
$ws.Cells.Item(87, 8).Value = 49999
$ws.Cells.Item(87, 10).Value = 49999
$ws.Cells.Item(87, 12).Value = 49999
$ws.Cells.Item(87, 14).Value = -52495
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 10).Value = 0
$ws.Cells.Item(88, 12).Value = 0
$ws.Cells.Item(89, 8).Value = 1983.875
$ws.Cells.Item(89, 9).Value = 1221.5
$ws.Cells.Item(89, 10).Value = 2746.25
$ws.Cells.Item(89, 11).Value = 6107.5
$ws.Cells.Item(89, 12).Value = 13731.25
$ws.Cells.Item(89, 13).Value = -491.5
$ws.Cells.Item(89, 14).Value = -24963.25
$ws.Cells.Item(90, 8).Value = 49999
$ws.Cells.Item(90, 10).Value = 49999
$ws.Cells.Item(90, 12).Value = 149997
$ws.Cells.Item(90, 14).Value = -162477
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 10).Value = 0
$ws.Cells.Item(91, 12).Value = 0
$ws.Cells.Item(93, 8).Value = 33448
$ws.Cells.Item(93, 10).Value = 33448
$ws.Cells.Item(93, 12).Value = 33448
$ws.Cells.Item(93, 14).Value = -37192
$ws.Cells.Item(94, 8).Value = 997.7647
$ws.Cells.Item(94, 9).Value = 845.25
$ws.Cells.Item(94, 10).Value = 1363.8
$ws.Cells.Item(94, 11).Value = 845.25
$ws.Cells.Item(94, 12).Value = 1363.8
$ws.Cells.Item(94, 13).Value = -394.25
$ws.Cells.Item(94, 14).Value = -2265.8
$ws.Cells.Item(95, 8).Value = 120000
$ws.Cells.Item(95, 10).Value = 120000
$ws.Cells.Item(95, 12).Value = 120000
$ws.Cells.Item(95, 14).Value = -125492
$ws.Cells.Item(96, 8).Value = 20139
$ws.Cells.Item(96, 9).Value = 20139
$ws.Cells.Item(96, 11).Value = 20139
$ws.Cells.Item(96, 13).Value = -17393
$ws.Cells.Item(97, 8).Value = 8185.778
$ws.Cells.Item(97, 9).Value = 8185.778
$ws.Cells.Item(97, 11).Value = 8185.778
$ws.Cells.Item(97, 13).Value = -7194.778
$ws.Cells.Item(98, 8).Value = 48333.332
$ws.Cells.Item(98, 10).Value = 48333.332
$ws.Cells.Item(98, 12).Value = 48333.332
$ws.Cells.Item(98, 14).Value = -54323.332
$ws.Cells.Item(99, 8).Value = 2022
$ws.Cells.Item(99, 9).Value = 706.9474
$ws.Cells.Item(99, 11).Value = 706.9474
$ws.Cells.Item(99, 13).Value = 791.0526
$ws.Cells.Item(101, 8).Value = 0
$ws.Cells.Item(101, 10).Value = 0
$ws.Cells.Item(101, 12).Value = 0
$ws.Cells.Item(104, 8).Value = 0
$ws.Cells.Item(104, 10).Value = 0
$ws.Cells.Item(104, 12).Value = 0
$ws.Cells.Item(105, 8).Value = 2562.9167
$ws.Cells.Item(105, 9).Value = 2520
$ws.Cells.Item(105, 10).Value = 2691.6667
$ws.Cells.Item(105, 11).Value = 2520
$ws.Cells.Item(105, 12).Value = 2691.6667
$ws.Cells.Item(105, 13).Value = -773
$ws.Cells.Item(105, 14).Value = -6185.6667
$ws.Cells.Item(106, 8).Value = 12500
$ws.Cells.Item(106, 10).Value = 12500
$ws.Cells.Item(106, 12).Value = 12500
$ws.Cells.Item(106, 14).Value = -15024
$ws.Cells.Item(134, 8).Value = 5553.25
$ws.Cells.Item(134, 9).Value = 5165.926
$ws.Cells.Item(134, 11).Value = 15497.778
$ws.Cells.Item(134, 13).Value = -12962.778
$ws.Cells.Item(88, 14).ClearContents()
$ws.Cells.Item(91, 14).ClearContents()
$ws.Cells.Item(101, 14).ClearContents()
$ws.Cells.Item(104, 14).ClearContents()

$ws = $wb.Worksheets("CRP")
$ws.Cells.Item(31, 8).Value = 2307.95
$ws.Cells.Item(31, 9).Value = 1752.8235
$ws.Cells.Item(31, 10).Value = 5453.6665
$ws.Cells.Item(31, 11).Value = 1752.8235
$ws.Cells.Item(31, 12).Value = 5453.6665
$ws.Cells.Item(31, 13).Value = -1457.8235
$ws.Cells.Item(31, 14).Value = -6043.6665
$ws.Cells.Item(34, 8).Value = 2307.95
$ws.Cells.Item(34, 9).Value = 1752.8235
$ws.Cells.Item(34, 10).Value = 5453.6665
$ws.Cells.Item(34, 11).Value = 1752.8235
$ws.Cells.Item(34, 12).Value = 5453.6665
$ws.Cells.Item(34, 13).Value = -1550.8235
$ws.Cells.Item(34, 14).Value = -5857.6665
$ws.Cells.Item(99, 8).Value = 8877.651
$ws.Cells.Item(99, 9).Value = 5562.5
$ws.Cells.Item(99, 10).Value = 13065.211
$ws.Cells.Item(99, 11).Value = 5562.5
$ws.Cells.Item(99, 12).Value = 13065.211
$ws.Cells.Item(99, 13).Value = -4064.5
$ws.Cells.Item(99, 14).Value = -16061.211
$ws.Cells.Item(126, 8).Value = 8877.651
$ws.Cells.Item(126, 9).Value = 5562.5
$ws.Cells.Item(126, 10).Value = 13065.211
$ws.Cells.Item(126, 11).Value = 16687.5
$ws.Cells.Item(126, 12).Value = 39195.633
$ws.Cells.Item(126, 13).Value = -14217.5
$ws.Cells.Item(126, 14).Value = -44135.633
$ws.Cells.Item(132, 8).Value = 5669.8047
$ws.Cells.Item(132, 9).Value = 3369.8064
$ws.Cells.Item(132, 11).Value = 10109.4192
$ws.Cells.Item(132, 13).Value = -7579.4192
$ws.Cells.Item(134, 8).Value = 2632.255
$ws.Cells.Item(134, 9).Value = 2426.625
$ws.Cells.Item(134, 11).Value = 7279.875
$ws.Cells.Item(134, 13).Value = -4744.875

$ws = $wb.Worksheets("CUL")
$ws.Cells.Item(38, 8).Value = 505.73914
$ws.Cells.Item(38, 9).Value = 369.55554
$ws.Cells.Item(38, 10).Value = 996
$ws.Cells.Item(38, 11).Value = 1108.66662
$ws.Cells.Item(38, 12).Value = 2988
$ws.Cells.Item(38, 13).Value = -761.66662
$ws.Cells.Item(38, 14).Value = -3682
$ws.Cells.Item(129, 8).Value = 1185.0667
$ws.Cells.Item(129, 9).Value = 616.3
$ws.Cells.Item(129, 10).Value = 2322.6
$ws.Cells.Item(129, 11).Value = 1848.9
$ws.Cells.Item(129, 12).Value = 6967.799999999999
$ws.Cells.Item(129, 13).Value = 3151.1
$ws.Cells.Item(129, 14).Value = -16967.8

$ws = $wb.Worksheets("GSM")
$ws.Cells.Item(46, 8).Value = 37970.75
$ws.Cells.Item(46, 9).Value = 4553.2
$ws.Cells.Item(46, 11).Value = 4553.2
$ws.Cells.Item(46, 13).Value = -4397.2
$ws.Cells.Item(126, 8).Value = 6598.933
$ws.Cells.Item(126, 9).Value = 5489.2
$ws.Cells.Item(126, 10).Value = 7153.8
$ws.Cells.Item(126, 11).Value = 16467.6
$ws.Cells.Item(126, 12).Value = 21461.4
$ws.Cells.Item(126, 13).Value = -13997.6
$ws.Cells.Item(126, 14).Value = -26401.4

$ws = $wb.Worksheets("LTW")
$ws.Cells.Item(29, 8).Value = 70016
$ws.Cells.Item(29, 10).Value = 0
$ws.Cells.Item(29, 12).Value = 0
$ws.Cells.Item(136, 8).Value = 17482.916
$ws.Cells.Item(136, 9).Value = 23079.4
$ws.Cells.Item(136, 11).Value = 69238.20000000001
$ws.Cells.Item(136, 13).Value = -66688.20000000001
$ws.Cells.Item(29, 14).ClearContents()

$ws = $wb.Worksheets("WVR")
$ws.Cells.Item(34, 8).Value = 26999.5
$ws.Cells.Item(34, 9).Value = 19999
$ws.Cells.Item(34, 10).Value = 34000
$ws.Cells.Item(34, 11).Value = 19999
$ws.Cells.Item(34, 12).Value = 34000
$ws.Cells.Item(34, 13).Value = -19796
$ws.Cells.Item(34, 14).Value = -34406
$ws.Cells.Item(70, 8).Value = 50000
$ws.Cells.Item(70, 10).Value = 50000
$ws.Cells.Item(70, 12).Value = 50000
$ws.Cells.Item(70, 14).Value = -50630
$ws.Cells.Item(73, 8).Value = 50000
$ws.Cells.Item(73, 10).Value = 50000
$ws.Cells.Item(73, 12).Value = 50000
$ws.Cells.Item(73, 14).Value = -52184
$ws.Cells.Item(126, 8).Value = 6582919.5
$ws.Cells.Item(126, 9).Value = 8337024.5
$ws.Cells.Item(126, 10).Value = 5025.75
$ws.Cells.Item(126, 11).Value = 25011073.5
$ws.Cells.Item(126, 12).Value = 15077.25
$ws.Cells.Item(126, 13).Value = -25008603.5
$ws.Cells.Item(126, 14).Value = -20017.25
$ws.Cells.Item(136, 8).Value = 2210.2341
$ws.Cells.Item(136, 9).Value = 2127.848
$ws.Cells.Item(136, 11).Value = 6383.544
$ws.Cells.Item(136, 13).Value = -3833.544
